# Add a new order line (row 7) to the sheet, matching the source data which
# stores every column - including the numeric-looking Quantity/Cost/Total
# columns - as text. Prefixing with a leading apostrophe forces Excel to
# store these as text (quote-prefixed) instead of auto-converting them to
# numbers, which matches how the rest of the sheet's rows are encoded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "SAB55321"
$ws.Range("B7").Value = "Container - Paper Clamshell (Bagel Box)"
$ws.Range("C7").Value = "'2"
$ws.Range("D7").Value = "'91.42"
$ws.Range("E7").Value = "'182.84"
